$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(1).Insert()
Write-Host "done"
